# Update LR-pair data with new TPM-based values and drop the two extra
# rows that the old (non-TPM) script used to produce.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the now-obsolete rows 7 and 6 (highest row first so indices
#     of the remaining rows are not disturbed while deleting). ---
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# --- Row 2 (ECs -> Ccl5/Ccr5 -> ECs) ---
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.021814
$ws.Range("H2").Value = 0.065442
$ws.Range("I2").Value = 0.1008129179549036
$ws.Range("J2").Value = 0.1008129179549036
$ws.Range("M2").Value = 0.001937666666666667
$ws.Range("N2").Value = 0.005813
$ws.Range("O2").Value = 0.0230007399171451
$ws.Range("P2").Value = 0.02300073991714511
$ws.Range("Q2").Value = 0.00004226826066666667
$ws.Range("R2").Value = 0.000380414346
$ws.Range("S2").Value = 0.002318771706169226
$ws.Range("T2").Value = 0.002318771706169226

# --- Row 3 (ECs -> Ccl5/Ccr5 -> FAPs) ---
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.021814
$ws.Range("H3").Value = 0.065442
$ws.Range("I3").Value = 0.1008129179549036
$ws.Range("J3").Value = 0.1008129179549036
$ws.Range("M3").Value = 0.082306
$ws.Range("O3").Value = 0.9769992600828549
$ws.Range("P3").Value = 0.976999260082855
$ws.Range("Q3").Value = 0.001795423084
$ws.Range("R3").Value = 0.016158807756
$ws.Range("S3").Value = 0.09849414624873438
$ws.Range("T3").Value = 0.0984941462487344

# --- Row 4: sending/target clusters swap (MuSCs -> Ccl5/Ccr5 -> ECs) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.194567
$ws.Range("H4").Value = 0.583701
$ws.Range("I4").Value = 0.8991870820450963
$ws.Range("J4").Value = 0.8991870820450963
$ws.Range("M4").Value = 0.001937666666666667
$ws.Range("N4").Value = 0.005813
$ws.Range("O4").Value = 0.0230007399171451
$ws.Range("P4").Value = 0.02300073991714511
$ws.Range("Q4").Value = 0.0003770059903333333
$ws.Range("R4").Value = 0.003393053913
$ws.Range("S4").Value = 0.02068196821097587
$ws.Range("T4").Value = 0.02068196821097588

# --- Row 5 (MuSCs -> Ccl5/Ccr5 -> FAPs) ---
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.194567
$ws.Range("H5").Value = 0.583701
$ws.Range("I5").Value = 0.8991870820450963
$ws.Range("J5").Value = 0.8991870820450963
$ws.Range("M5").Value = 0.082306
$ws.Range("N5").Value = 0.246918
$ws.Range("O5").Value = 0.9769992600828549
$ws.Range("P5").Value = 0.976999260082855
$ws.Range("Q5").Value = 0.016014031502
$ws.Range("R5").Value = 0.144126283518
$ws.Range("S5").Value = 0.8785051138341204
$ws.Range("T5").Value = 0.8785051138341206

Write-Output "Applied NATMI TPM update"
